# Refresh Ligand / Receptor / Edge expression + specificity metrics
# (columns G:T) for the new TPM-based NATMI run -- rows 2-26 of the
# Col4a1 -> Itgb1 ligand-receptor pairs sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowValues([int]$row, [object[]]$vals) {
    $arr = New-Object "object[,]" 1,$vals.Length
    for ($i = 0; $i -lt $vals.Length; $i++) { $arr[0,$i] = $vals[$i] }
    $ws.Range("G$row`:T$row").Value = $arr
}

Set-RowValues 2 @(421.7027586666666, 1265.108276, 0.8230205889389596, 0.8230205889389596, 3, 1, 159.4836373333333, 478.450912, 0.2983285084902258, 0.2983285084902258, 67254.68982566085, 605292.2084309477, 0.245530504754907, 0.245530504754907)
Set-RowValues 3 @(421.7027586666666, 1265.108276, 0.8230205889389596, 0.8230205889389596, 3, 1, 172.558497, 517.675491, 0.3227862111630279, 0.3227862111630279, 72768.39421627371, 654915.5479464635, 0.2656596976127706, 0.2656596976127706)
Set-RowValues 4 @(421.7027586666666, 1265.108276, 0.8230205889389596, 0.8230205889389596, 3, 1, 74.38770566666666, 223.163117, 0.1391489036280481, 0.1391489036280482, 31369.50069051736, 282325.5062146563, 0.1145224126141667, 0.1145224126141667)
Set-RowValues 5 @(421.7027586666666, 1265.108276, 0.8230205889389596, 0.8230205889389596, 3, 1, 58.41461433333333, 175.243843, 0.1092697975759847, 0.1092697975759848, 24633.60401081607, 221702.4360973446, 0.08993129315422786, 0.08993129315422788)
Set-RowValues 6 @(421.7027586666666, 1265.108276, 0.8230205889389596, 0.8230205889389596, 3, 1, 69.746216, 209.238648, 0.1304665791427133, 0.1304665791427133, 29412.1716937612, 264709.5452438509, 0.1073766808028873, 0.1073766808028873)
Set-RowValues 7 @(69.70494733333334, 209.114842, 0.136040387754698, 0.136040387754698, 3, 1, 159.4836373333333, 478.450912, 0.2983285084902258, 0.2983285084902258, 11116.79854084843, 100051.1868676359, 0.04058472597329103, 0.04058472597329103)
Set-RowValues 8 @(69.70494733333334, 209.114842, 0.136040387754698, 0.136040387754698, 3, 1, 172.558497, 517.675491, 0.3227862111630279, 0.3227862111630279, 12028.18094530416, 108253.6285077374, 0.04391196132848815, 0.04391196132848815)
Set-RowValues 9 @(69.70494733333334, 209.114842, 0.136040387754698, 0.136040387754698, 3, 1, 74.38770566666666, 223.163117, 0.1391489036280481, 0.1391489036280482, 5185.191105742501, 46666.71995168251, 0.01892987080520077, 0.01892987080520078)
Set-RowValues 10 @(69.70494733333334, 209.114842, 0.136040387754698, 0.136040387754698, 3, 1, 58.41461433333333, 175.243843, 0.1092697975759847, 0.1092697975759848, 4071.787615601979, 36646.0885404178, 0.01486510563211433, 0.01486510563211433)
Set-RowValues 11 @(69.70494733333334, 209.114842, 0.136040387754698, 0.136040387754698, 3, 1, 69.746216, 209.238648, 0.1304665791427133, 0.1304665791427133, 4861.656312979291, 43754.90681681362, 0.01774872401560372, 0.01774872401560372)
Set-RowValues 12 @(1.088159666666667, 3.264479, 0.002123718167154624, 0.002123718167154624, 3, 1, 159.4836373333333, 478.450912, 0.2983285084902258, 0.2983285084902258, 173.5436616394276, 1561.892954754848, 0.000633565673260835, 0.000633565673260835)
Set-RowValues 13 @(1.088159666666667, 3.264479, 0.002123718167154624, 0.002123718167154624, 3, 1, 172.558497, 517.675491, 0.3227862111630279, 0.3227862111630279, 187.771196576021, 1689.940769184189, 0.0006855069407539311, 0.0006855069407539311)
Set-RowValues 14 @(1.088159666666667, 3.264479, 0.002123718167154624, 0.002123718167154624, 3, 1, 74.38770566666666, 223.163117, 0.1391489036280481, 0.1391489036280482, 80.9457010023381, 728.5113090210431, 0.0002955130545745338, 0.0002955130545745339)
Set-RowValues 15 @(1.088159666666667, 3.264479, 0.002123718167154624, 0.002123718167154624, 3, 1, 58.41461433333333, 175.243843, 0.1092697975759847, 0.1092697975759848, 63.56442726142188, 572.0798453527971, 0.0002320582542334271, 0.0002320582542334271)
Set-RowValues 16 @(1.088159666666667, 3.264479, 0.002123718167154624, 0.002123718167154624, 3, 1, 69.746216, 209.238648, 0.1304665791427133, 0.1304665791427133, 75.89501915382134, 683.0551723843921, 0.0002770742443318968, 0.0002770742443318969)
Set-RowValues 17 @(19.43698366666667, 58.310951, 0.03793439197579861, 0.03793439197579861, 3, 1, 159.4836373333333, 478.450912, 0.2983285084902258, 0.2983285084902258, 3099.88085394859, 27898.92768553731, 0.01131691057862359, 0.01131691057862359)
Set-RowValues 18 @(19.43698366666667, 58.310951, 0.03793439197579861, 0.03793439197579861, 3, 1, 172.558497, 517.675491, 0.3227862111630279, 0.3227862111630279, 3354.016687733549, 30186.15018960194, 0.0122446986586412, 0.0122446986586412)
Set-RowValues 19 @(19.43698366666667, 58.310951, 0.03793439197579861, 0.03793439197579861, 3, 1, 74.38770566666666, 223.163117, 0.1391489036280481, 0.1391489036280482, 1445.872620043807, 13012.85358039427, 0.005278529053229003, 0.005278529053229005)
Set-RowValues 20 @(19.43698366666667, 58.310951, 0.03793439197579861, 0.03793439197579861, 3, 1, 58.41461433333333, 175.243843, 0.1092697975759847, 0.1092697975759848, 1135.403904691633, 10218.63514222469, 0.004145083332363575, 0.004145083332363575)
Set-RowValues 21 @(19.43698366666667, 58.310951, 0.03793439197579861, 0.03793439197579861, 3, 1, 69.746216, 209.238648, 0.1304665791427133, 0.1304665791427133, 1355.656061203805, 12200.90455083425, 0.004949170352941239, 0.00494917035294124)
Set-RowValues 22 @(0.451366, 1.354098, 0.0008809131633892397, 0.0008809131633892397, 3, 1, 159.4836373333333, 478.450912, 0.2983285084902258, 0.2983285084902258, 71.98549144859734, 647.869423037376, 0.0002628015101433185, 0.0002628015101433185)
Set-RowValues 23 @(0.451366, 1.354098, 0.0008809131633892397, 0.0008809131633892397, 3, 1, 172.558497, 517.675491, 0.3227862111630279, 0.3227862111630279, 77.887038556902, 700.983347012118, 0.0002843466223740501, 0.0002843466223740501)
Set-RowValues 24 @(0.451366, 1.354098, 0.0008809131633892397, 0.0008809131633892397, 3, 1, 74.38770566666666, 223.163117, 0.1391489036280481, 0.1391489036280482, 33.57608115594066, 302.184730403466, 0.0001225781008771283, 0.0001225781008771284)
Set-RowValues 25 @(0.451366, 1.354098, 0.0008809131633892397, 0.0008809131633892397, 3, 1, 58.41461433333333, 175.243843, 0.1092697975759847, 0.1092697975759848, 26.36637081317933, 237.297337318614, 0.00009625720304556261, 0.00009625720304556262)
Set-RowValues 26 @(0.451366, 1.354098, 0.0008809131633892397, 0.0008809131633892397, 3, 1, 69.746216, 209.238648, 0.1304665791427133, 0.1304665791427133, 31.481070531056, 283.329634779504, 0.0001149297269491802, 0.0001149297269491802)

Write-Host "Updated rows 2-26, columns G:T with new TPM-derived values"
